$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ErrorDetails")

# Path string reused across rows
$path = "C:\Users\Siddharth Shinde\Desktop\Springmaven\myQfcProject\myQfcProject\mycommon\src\main\resources\success.xlsx"

# Update row 2 (existing row) with new values
$ws.Range("A2").Value = "sidd"
$ws.Range("B2").Value = "Invalid numeric cell type"
$ws.Range("C2").Value = 6.0
$ws.Range("D2").Value = "UserID"
$ws.Range("E2").Value = $path

# Add row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "12.0"
$ws.Range("B3").Value = "Invalid string cell type"
$ws.Range("C3").Value = 8.0
$ws.Range("D3").Value = "UserName"
$ws.Range("E3").Value = $path

# Add row 4 (same as row 2)
$ws.Range("A4").Value = "sidd"
$ws.Range("B4").Value = "Invalid numeric cell type"
$ws.Range("C4").Value = 6.0
$ws.Range("D4").Value = "UserID"
$ws.Range("E4").Value = $path

# Add row 5 (same as row 3)
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "12.0"
$ws.Range("B5").Value = "Invalid string cell type"
$ws.Range("C5").Value = 8.0
$ws.Range("D5").Value = "UserName"
$ws.Range("E5").Value = $path
